$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 3: "T2TRG/IETF (20m)" -> "T2TRG/IETF (20m) – Carsten Bormann"
$para3 = $tr.Paragraphs(3, 1)
$run3 = $para3.Runs(1, 1)
$run3.Text = "T2TRG/IETF (20m) – Carsten Bormann"

# Paragraph 4: "ECHONET (20m)" -> "ECHONET (20m) – Tetsushi Matsuda" (3 runs)
$tr = $shape.TextFrame.TextRange
$para4 = $tr.Paragraphs(4, 1)
$run4a = $para4.Runs(1, 1)
$run4a.Text = "ECHONET (20m) – "
$run4b = $run4a.InsertAfter("Tetsushi")
$run4c = $run4b.InsertAfter(" Matsuda")

# Paragraph 5: "ITU-T (20m)" -> "ITU-T (20m) – Gyu Myoung Lee" (5 runs)
$tr = $shape.TextFrame.TextRange
$para5 = $tr.Paragraphs(5, 1)
$run5a = $para5.Runs(1, 1)
$run5a.Text = "ITU-T (20m) – "
$run5b = $run5a.InsertAfter("Gyu")
$run5c = $run5b.InsertAfter(" ")
$run5d = $run5c.InsertAfter("Myoung")
$run5e = $run5d.InsertAfter(" Lee")

# Paragraph 8: "Eclipse Ditto (Thomas Jäckle, 20m)" -> "Eclipse Ditto (20m) - Thomas Jäckle"
# (the trailing ", 20m)" run is removed so the paragraph now ends right after "Jäckle")
$tr = $shape.TextFrame.TextRange
$para8 = $tr.Paragraphs(8, 1)
$run8a = $para8.Runs(1, 1)
$run8a.Text = "Eclipse Ditto (20m) - Thomas "
# run 2 ("Jäckle") is left untouched so its formatting (incl. err="1") is preserved
$run8c = $para8.Runs(3, 1)
$run8c.Text = ""

# Paragraph 9: "ISO TC184/SC4 (about IEC CDD) (20m)" ->
#   "ISO TC184/SC4 (about IEC CDD) (20m) – Hiroshi Murayama/Yoshiaki Sonoda" (2 runs)
$tr = $shape.TextFrame.TextRange
$para9 = $tr.Paragraphs(9, 1)
$run9a = $para9.Runs(1, 1)
$run9a.Text = "ISO TC184/SC4 (about IEC CDD) (20m) – Hiroshi Murayama/Yoshiaki "
$run9b = $run9a.InsertAfter("Sonoda")
